$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A10").Value = "admin1"
$ws.Range("B10").Value = "123456a!"
$ws.Range("C10").Value = 309308906
